$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "75.777.37"
$ws.Range("E2").Value = "  +8.91%  "
$ws.Range("D3").Value = "2.690.75"
$ws.Range("E3").Value = "  +10.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.93"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +13.56%  "
$ws.Range("E6").Value = "  +4.50%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +5.37%  "
$ws.Range("E9").Value = "  +15.09%  "
$ws.Range("D10").Value = "2.689.90"
$ws.Range("E10").Value = "  +10.95%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("E12").Value = "  +7.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.75"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "75.704.80"
$ws.Range("E14").Value = "  +8.99%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.188.04"
$ws.Range("E15").Value = "  +11.03%  "
$ws.Range("E16").Value = "  +6.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.58"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  +10.86%  "
$ws.Range("D18").Value = "2.685.71"
$ws.Range("E18").Value = "  +10.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.37"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +31.14%  "
$ws.Range("E20").Value = "  +11.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.31"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +9.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.31"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  +17.58%  "
$ws.Range("E23").Value = "  +5.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.27"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.33"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +6.16%  "
$ws.Range("E27").Value = "  +9.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.44"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  +11.10%  "
$ws.Range("D29").Value = "2.826.28"
$ws.Range("E29").Value = "  +10.55%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "0.0₃0953"
$ws.Range("E31").Value = "  +12.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "522.53"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  +14.93%  "
$ws.Range("E33").Value = "  +13.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.77"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  +5.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.77"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +9.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.74"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("E38").Value = "  +8.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.31"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  +6.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.39"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.04"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +14.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "170.92"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +26.40%  "
$ws.Range("E44").Value = "  +12.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.333"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +9.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  +10.48%  "
$ws.Range("E47").Value = "  +14.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.41"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.67"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +8.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.540"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +10.38%  "
